$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 1621.75
$ws.Cells.Item(38, 9).Value = 46.1
$ws.Cells.Item(38, 10).Value = 9500
$ws.Cells.Item(38, 11).Value = 138.3
$ws.Cells.Item(38, 12).Value = 28500
$ws.Cells.Item(38, 13).Value = 233.7
$ws.Cells.Item(38, 14).Value = -29244
$ws.Cells.Item(39, 8).Value = 147.8125
$ws.Cells.Item(39, 9).Value = 91
$ws.Cells.Item(39, 11).Value = 273
$ws.Cells.Item(39, 13).Value = 23
$ws.Cells.Item(55, 8).Value = 1259.2727
$ws.Cells.Item(55, 9).Value = 595.1667
$ws.Cells.Item(55, 10).Value = 2056.2
$ws.Cells.Item(55, 11).Value = 595.1667
$ws.Cells.Item(55, 12).Value = 2056.2
$ws.Cells.Item(55, 13).Value = -381.1667
$ws.Cells.Item(55, 14).Value = -2484.2
$ws.Cells.Item(88, 8).Value = 2333
$ws.Cells.Item(88, 10).Value = 1999.5
$ws.Cells.Item(88, 12).Value = 1999.5
$ws.Cells.Item(88, 14).Value = -2811.5
$ws.Cells.Item(91, 8).Value = 2333
$ws.Cells.Item(91, 10).Value = 1999.5
$ws.Cells.Item(91, 12).Value = 1999.5
$ws.Cells.Item(91, 14).Value = -4807.5
$ws.Cells.Item(106, 8).Value = 6998.5
$ws.Cells.Item(106, 9).Value = 6998.5
$ws.Cells.Item(106, 11).Value = 6998.5
$ws.Cells.Item(106, 13).Value = -6367.5
$ws.Cells.Item(132, 8).Value = 34451.375
$ws.Cells.Item(132, 9).Value = 42122.2
$ws.Cells.Item(132, 11).Value = 126366.6
$ws.Cells.Item(132, 13).Value = -123836.6
$ws.Cells.Item(135, 8).Value = 855.5833
$ws.Cells.Item(135, 9).Value = 698.2
$ws.Cells.Item(135, 10).Value = 1642.5
$ws.Cells.Item(135, 11).Value = 6283.8
$ws.Cells.Item(135, 12).Value = 14782.5
$ws.Cells.Item(135, 13).Value = -3748.8
$ws.Cells.Item(135, 14).Value = -19852.5
$ws.Cells.Item(137, 8).Value = 2092.1667
$ws.Cells.Item(137, 9).Value = 992.1818
$ws.Cells.Item(137, 11).Value = 2976.5454
$ws.Cells.Item(137, 13).Value = -426.5454
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2437.3635
$ws.Cells.Item(2, 9).Value = 2762.2
$ws.Cells.Item(2, 10).Value = 2166.6667
$ws.Cells.Item(2, 11).Value = 2762.2
$ws.Cells.Item(2, 12).Value = 2166.6667
$ws.Cells.Item(2, 13).Value = -2649.2
$ws.Cells.Item(2, 14).Value = -2392.6667
$ws.Cells.Item(45, 8).Value = 2982.2144
$ws.Cells.Item(45, 9).Value = 2451.4443
$ws.Cells.Item(45, 10).Value = 3937.6
$ws.Cells.Item(45, 11).Value = 2451.4443
$ws.Cells.Item(45, 12).Value = 3937.6
$ws.Cells.Item(45, 13).Value = -2074.4443
$ws.Cells.Item(45, 14).Value = -4691.6
$ws.Cells.Item(102, 8).Value = 12506788
$ws.Cells.Item(102, 9).Value = 41669964
$ws.Cells.Item(102, 11).Value = 41669964
$ws.Cells.Item(102, 13).Value = -41668342
$ws.Cells.Item(110, 8).Value = 111113430
$ws.Cells.Item(110, 9).Value = 166668670
$ws.Cells.Item(110, 11).Value = 166668670
$ws.Cells.Item(110, 13).Value = -166666625
$ws.Cells.Item(116, 8).Value = 2437.3635
$ws.Cells.Item(116, 9).Value = 2762.2
$ws.Cells.Item(116, 10).Value = 2166.6667
$ws.Cells.Item(116, 11).Value = 2762.2
$ws.Cells.Item(116, 12).Value = 2166.6667
$ws.Cells.Item(116, 13).Value = -468.1999999999998
$ws.Cells.Item(116, 14).Value = -6754.6667
$ws.Cells.Item(132, 8).Value = 1649
$ws.Cells.Item(132, 9).Value = 1313.1428
$ws.Cells.Item(132, 11).Value = 3939.4284
$ws.Cells.Item(132, 13).Value = -1409.4284
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2437.3635
$ws.Cells.Item(3, 9).Value = 2762.2
$ws.Cells.Item(3, 10).Value = 2166.6667
$ws.Cells.Item(3, 11).Value = 2762.2
$ws.Cells.Item(3, 12).Value = 2166.6667
$ws.Cells.Item(3, 13).Value = -2648.2
$ws.Cells.Item(3, 14).Value = -2394.6667
$ws.Cells.Item(20, 8).Value = 3956
$ws.Cells.Item(20, 9).Value = 3723.5
$ws.Cells.Item(20, 11).Value = 3723.5
$ws.Cells.Item(20, 13).Value = -3476.5
$ws.Cells.Item(68, 8).Value = 35000
$ws.Cells.Item(68, 10).Value = 35000
$ws.Cells.Item(68, 12).Value = 35000
$ws.Cells.Item(68, 14).Value = -36622
$ws.Cells.Item(71, 8).Value = 35000
$ws.Cells.Item(71, 10).Value = 35000
$ws.Cells.Item(71, 12).Value = 105000
$ws.Cells.Item(71, 14).Value = -113112
$ws.Cells.Item(86, 8).Value = 6139.8184
$ws.Cells.Item(86, 9).Value = 1907.8
$ws.Cells.Item(86, 11).Value = 1907.8
$ws.Cells.Item(86, 13).Value = -784.8
$ws.Cells.Item(89, 8).Value = 6139.8184
$ws.Cells.Item(89, 9).Value = 1907.8
$ws.Cells.Item(89, 11).Value = 9539
$ws.Cells.Item(89, 13).Value = -3923
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(29, 8).Value = 2500
$ws.Cells.Item(29, 9).Value = 2500
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 2500
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = -2207
$ws.Cells.Item(29, 14).ClearContents()
$ws.Cells.Item(105, 8).Value = 2955.625
$ws.Cells.Item(105, 9).Value = 1326.8
$ws.Cells.Item(105, 11).Value = 1326.8
$ws.Cells.Item(105, 13).Value = 420.2
$ws.Cells.Item(107, 8).Value = 1610.8462
$ws.Cells.Item(107, 9).Value = 583.5
$ws.Cells.Item(107, 11).Value = 583.5
$ws.Cells.Item(107, 13).Value = 1336.5
$ws.Cells.Item(138, 8).Value = 97500
$ws.Cells.Item(138, 10).Value = 97500
$ws.Cells.Item(138, 12).Value = 97500
$ws.Cells.Item(138, 14).Value = -107780
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 5
$ws.Cells.Item(7, 10).Value = 5
$ws.Cells.Item(7, 12).Value = 15
$ws.Cells.Item(7, 14).Value = -239
$ws.Cells.Item(16, 8).Value = 566.8333
$ws.Cells.Item(16, 9).Value = 100.25
$ws.Cells.Item(16, 10).Value = 1500
$ws.Cells.Item(16, 11).Value = 300.75
$ws.Cells.Item(16, 12).Value = 4500
$ws.Cells.Item(16, 13).Value = -127.75
$ws.Cells.Item(16, 14).Value = -4846
$ws.Cells.Item(80, 8).Value = 4245.625
$ws.Cells.Item(80, 9).Value = 3977.5789
$ws.Cells.Item(80, 11).Value = 11932.7367
$ws.Cells.Item(80, 13).Value = -10996.7367
$ws.Cells.Item(83, 8).Value = 4245.625
$ws.Cells.Item(83, 9).Value = 3977.5789
$ws.Cells.Item(83, 11).Value = 35798.2101
$ws.Cells.Item(83, 13).Value = -31118.2101
$ws.Cells.Item(92, 8).Value = 1662.4546
$ws.Cells.Item(92, 10).Value = 2164.5
$ws.Cells.Item(92, 12).Value = 6493.5
$ws.Cells.Item(92, 14).Value = -8989.5
$ws.Cells.Item(129, 8).Value = 1738
$ws.Cells.Item(129, 9).Value = 1566.6666
$ws.Cells.Item(129, 10).Value = 1995
$ws.Cells.Item(129, 11).Value = 4699.9998
$ws.Cells.Item(129, 12).Value = 5985
$ws.Cells.Item(129, 13).Value = 300.0002000000004
$ws.Cells.Item(129, 14).Value = -15985
$ws.Cells.Item(131, 8).Value = 2412.3333
$ws.Cells.Item(131, 10).Value = 2559.9
$ws.Cells.Item(131, 12).Value = 7679.700000000001
$ws.Cells.Item(131, 14).Value = -17759.7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3565
$ws.Cells.Item(81, 9).Value = 1416.2
$ws.Cells.Item(81, 10).Value = 6251
$ws.Cells.Item(81, 11).Value = 2832.4
$ws.Cells.Item(81, 12).Value = 12502
$ws.Cells.Item(81, 13).Value = -1771.4
$ws.Cells.Item(81, 14).Value = -14624
$ws.Cells.Item(84, 8).Value = 3565
$ws.Cells.Item(84, 9).Value = 1416.2
$ws.Cells.Item(84, 10).Value = 6251
$ws.Cells.Item(84, 11).Value = 14162
$ws.Cells.Item(84, 12).Value = 62510
$ws.Cells.Item(84, 13).Value = -8858
$ws.Cells.Item(84, 14).Value = -73118
$ws.Cells.Item(107, 8).Value = 83334090
$ws.Cells.Item(107, 9).Value = 83334090
$ws.Cells.Item(107, 11).Value = 250002270
$ws.Cells.Item(107, 13).Value = -250000350
$ws.Cells.Item(113, 8).Value = 1153.2727
$ws.Cells.Item(113, 9).Value = 981.3333
$ws.Cells.Item(113, 11).Value = 2943.9999
$ws.Cells.Item(113, 13).Value = -773.9998999999998
$ws.Cells.Item(129, 8).Value = 110000
$ws.Cells.Item(129, 10).Value = 110000
$ws.Cells.Item(129, 12).Value = 110000
$ws.Cells.Item(129, 14).Value = -120000
$ws.Cells.Item(132, 8).Value = 1981.2273
$ws.Cells.Item(132, 9).Value = 1879.35
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 5638.049999999999
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = -3108.049999999999
$ws.Cells.Item(132, 14).Value = -14060
